$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-07-22 Monday" "2024-07-23 Tuesday"

Replace-Text "675×5=" "607×4="
Replace-Text "414×5=" "519×7="
Replace-Text "668×7=" "287×3="
Replace-Text "744×8=" "291×3="
Replace-Text "799×6=" "261×2="

Replace-Text "210×2=" "180×8="
Replace-Text "350×9=" "130×8="
Replace-Text "936×6=" "876×4="
Replace-Text "156×6=" "982×3="
Replace-Text "659×3=" "321×6="

Replace-Text "493×4=" "134×9="
Replace-Text "706×6=" "707×4="
Replace-Text "231×8=" "872×6="
Replace-Text "236×4=" "523×5="
Replace-Text "279×3=" "694×9="

Replace-Text "339×2=" "461×3="
Replace-Text "716×4=" "714×4="
Replace-Text "877×3=" "911×8="
Replace-Text "630×7=" "626×8="
Replace-Text "807×3=" "813×3="

Replace-Text "747×9=" "394×9="
Replace-Text "333×2=" "552×8="
Replace-Text "760×2=" "952×4="
Replace-Text "377×9=" "312×9="
Replace-Text "593×5=" "715×6="
